$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Style B1: bold font, thin box border, center horizontal / top vertical alignment
$cell1 = $ws.Range("B1")
$cell1.Font.Bold = $true
$cell1.HorizontalAlignment = -4108  # xlCenter
$cell1.VerticalAlignment = -4160    # xlTop
$cell1.Borders.LineStyle = 1        # xlContinuous
$cell1.Borders.Weight = 2           # xlThin

# Apply the same formatting to A2 by copying B1's format (avoids creating
# redundant intermediate cell-style entries)
$cell1.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
